# Updates the cryptocurrency listing (prices / 1h volume % / two row swaps)
# to match the latest scrape, preserving each cell's original text-string
# type (these columns are stored as text, not numbers, in the workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $ref, $value) {
    $cell = $ws.Range($ref)
    # Force a text number-format before assigning so Excel's COM layer does
    # not auto-coerce numeric-looking strings (e.g. "0.999", "586.15")
    # into real numbers, which would change the cell's stored type.
    $cell.NumberFormat = "@"
    $cell.Value = $value
    # Reset back to the default "Normal" style so no stray per-cell style
    # (beyond the text number format momentarily needed above) is left
    # behind on the saved cell.
    $cell.Style = "Normal"
}

Set-TextCell $ws 'D2' '61.670.54'
Set-TextCell $ws 'E2' '  -1.56%  '
Set-TextCell $ws 'D3' '2.901.94'
Set-TextCell $ws 'E3' '  -2.05%  '
Set-TextCell $ws 'D4' '0.999'
Set-TextCell $ws 'E4' '  -0.02%  '
Set-TextCell $ws 'D5' '586.15'
Set-TextCell $ws 'E5' '  -1.57%  '
Set-TextCell $ws 'D6' '146.01'
Set-TextCell $ws 'E6' '  +0.51%  '
Set-TextCell $ws 'E7' '  -0.02%  '
Set-TextCell $ws 'E8' '  +0.82%  '
Set-TextCell $ws 'D9' '2.901.47'
Set-TextCell $ws 'E9' '  -2.01%  '
Set-TextCell $ws 'E10' '  -6.11%  '
Set-TextCell $ws 'D11' '0.148'
Set-TextCell $ws 'E11' '  +3.61%  '
Set-TextCell $ws 'E12' '  -2.98%  '
Set-TextCell $ws 'D13' '0.0000236'
Set-TextCell $ws 'E13' '  +1.00%  '
Set-TextCell $ws 'E14' '  -1.87%  '
Set-TextCell $ws 'E15' '  -0.81%  '
Set-TextCell $ws 'D16' '3.383.84'
Set-TextCell $ws 'E16' '  -2.00%  '
Set-TextCell $ws 'D17' '61.671.55'
Set-TextCell $ws 'E17' '  -1.37%  '
Set-TextCell $ws 'E18' '  -1.95%  '
Set-TextCell $ws 'D19' '2.905.15'
Set-TextCell $ws 'E19' '  -1.90%  '
Set-TextCell $ws 'D20' '435.07'
Set-TextCell $ws 'E20' '  -1.37%  '
Set-TextCell $ws 'E21' '  -0.58%  '
Set-TextCell $ws 'E22' '  -2.06%  '
Set-TextCell $ws 'E23' '  -2.61%  '
Set-TextCell $ws 'D24' '81.07'
Set-TextCell $ws 'E24' '  -0.83%  '
Set-TextCell $ws 'D25' '11.94'
Set-TextCell $ws 'E25' '  +0.01%  '
Set-TextCell $ws 'D26' '10.23'
Set-TextCell $ws 'E26' '  -7.89%  '
Set-TextCell $ws 'E27' '  +0.00%  '
Set-TextCell $ws 'E28' '  -2.96%  '
Set-TextCell $ws 'B29' 'PEPE'
Set-TextCell $ws 'C29' 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextCell $ws 'D29' '0.0000104'
Set-TextCell $ws 'E29' '  +18.74%  '
Set-TextCell $ws 'B30' 'NEARProtocol'
Set-TextCell $ws 'C30' 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextCell $ws 'D30' '7.19'
Set-TextCell $ws 'E30' '  +1.59%  '
Set-TextCell $ws 'E31' '  -2.35%  '
Set-TextCell $ws 'E32' '  -1.23%  '
Set-TextCell $ws 'D33' '0.109'
Set-TextCell $ws 'E33' '  +0.64%  '
Set-TextCell $ws 'E34' '  +0.00%  '
Set-TextCell $ws 'D35' '25.88'
Set-TextCell $ws 'E35' '  -2.14%  '
Set-TextCell $ws 'E36' '  -1.90%  '
Set-TextCell $ws 'E37' '  -2.38%  '
Set-TextCell $ws 'E38' '  +2.69%  '
Set-TextCell $ws 'E39' '  -1.04%  '
Set-TextCell $ws 'E40' '  -3.16%  '
Set-TextCell $ws 'B41' 'Cosmos'
Set-TextCell $ws 'C41' 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextCell $ws 'D41' '8.35'
Set-TextCell $ws 'E41' '  -2.73%  '
Set-TextCell $ws 'B42' 'Kaspa'
Set-TextCell $ws 'C42' 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextCell $ws 'D42' '0.116'
Set-TextCell $ws 'E42' '  -1.26%  '
Set-TextCell $ws 'E43' '  -3.77%  '
Set-TextCell $ws 'D44' '38.81'
Set-TextCell $ws 'E44' '  -1.31%  '
Set-TextCell $ws 'D45' '2.690.42'
Set-TextCell $ws 'E45' '  -1.04%  '
Set-TextCell $ws 'D46' '133.70'
Set-TextCell $ws 'E46' '  -0.84%  '
Set-TextCell $ws 'D47' '0.0335'
Set-TextCell $ws 'E47' '  -1.91%  '
Set-TextCell $ws 'B48' 'Bittensor'
Set-TextCell $ws 'C48' 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextCell $ws 'D48' '341.95'
Set-TextCell $ws 'E48' '  -6.51%  '
Set-TextCell $ws 'B49' 'USDe'
Set-TextCell $ws 'C49' 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextCell $ws 'D49' '1.00'
Set-TextCell $ws 'E49' '  +0.01%  '
Set-TextCell $ws 'E50' '  -1.78%  '
Set-TextCell $ws 'D51' '22.25'
Set-TextCell $ws 'E51' '  -3.50%  '
